$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; existing B (Valor data) shifts to C
$ws.Columns("B:B").Insert()

# Header row
$ws.Range("B1").Value = "Variável"
$ws.Range("C1").Value = "Valor"
$ws.Range("D1").Value = "Colocação"

# New "Variável" column text for each data row
$ws.Range("B2:B10").Value = "Diferença 2021-2012"

# Ranking column "Colocação" for rows 2-8 (individual states; the two
# aggregate rows - Nordeste/9 and Brasil/10 - are left blank, as in the source)
$rankings = "1º", "2º", "3º", "4º", "5º", "6º", "13º"
for ($i = 0; $i -lt $rankings.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 4).Value = $rankings[$i]
}
